$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: region names (only rows 3-7 change)
$ws.Range("A3").Value = "Piauí"
$ws.Range("A4").Value = "Maranhão"
$ws.Range("A5").Value = "Tocantins"
$ws.Range("A6").Value = "Amazonas"
$ws.Range("A7").Value = "Paraíba"

# Column B: variable label text, all rows 2-10 change
$ws.Range("B2:B10").Value = "Diferença 2022-2000"

# Column C: values, all rows 2-10 change
$ws.Range("C2").Value = 0.263
$ws.Range("C3").Value = 0.258
$ws.Range("C4").Value = 0.252
$ws.Range("C5").Value = 0.252
$ws.Range("C6").Value = 0.249
$ws.Range("C7").Value = 0.232
$ws.Range("C8").Value = 0.223
$ws.Range("C9").Value = 0.196
$ws.Range("C10").Value = 0.1860000000000001

# Column D: ranking, only row 8 changes
$ws.Range("D8").Value = "10º"
